# "Added script for cash payments"
#
# - test_suite: enable the PayThruCash and MarkPaymentAsPaid test cases
#   (runmode N -> Y), and disable the PayThruAliPay test case (runmode Y -> N).
# - BillingAddress: flip the runmode flag to Y and drop the no-longer-used
#   "Reference Number" column (column N).
# - Leave the workbook with test_suite as the active sheet/selection instead
#   of OrderDetails / BillingAddress.

$wb = $excel.ActiveWorkbook

# --- test_suite ---------------------------------------------------------
$wsTest = $wb.Worksheets.Item("test_suite")
$wsTest.Range("B2").Value = "Y"   # PayThruCash
$wsTest.Range("B3").Value = "Y"   # MarkPaymentAsPaid
$wsTest.Range("B8").Value = "N"   # PayThruAliPay

# --- BillingAddress ------------------------------------------------------
$wsBilling = $wb.Worksheets.Item("BillingAddress")
$wsBilling.Range("M2").Value = "Y"     # runmode
$wsBilling.Columns.Item(14).Delete()   # remove "Reference Number" column (N)
$wsBilling.Range("F19").Select()

# --- Active sheet / selection housekeeping -------------------------------
$wsTest.Activate()
$wsTest.Range("B4").Select()
